# Weekly CompStat update: new crime data collected for week of 6/26/2023 - 7/2/2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlPasteFormats = -4122

# --- Header text: volume/issue number and reporting week date range ---
$ws.Range("A8").Value = "Volume 30   Number  26"
$ws.Range("C9").Value = "Report Covering the Week  6/26/2023  Through  7/2/2023"

# --- Plain numeric value updates (cell already numeric, style unchanged) ---
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("N15").Value = -64.285714285714
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = -20
$ws.Range("J16").Value = 55
$ws.Range("K16").Value = 5.454545454545
$ws.Range("L16").Value = 16
$ws.Range("M16").Value = -15.942028985507
$ws.Range("N16").Value = -85.204081632653
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = 33.333333333333
$ws.Range("I17").Value = 96
$ws.Range("J17").Value = 67
$ws.Range("K17").Value = 43.283582089552
$ws.Range("L17").Value = 104.255319148936
$ws.Range("M17").Value = 433.333333333333
$ws.Range("N17").Value = -17.948717948717
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -18.181818181818
$ws.Range("J18").Value = 60
$ws.Range("K18").Value = 60
$ws.Range("L18").Value = 50
$ws.Range("M18").Value = 9.090909090909
$ws.Range("N18").Value = -74.262734584450
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 17
$ws.Range("E19").Value = -35.294117647058
$ws.Range("F19").Value = 52
$ws.Range("G19").Value = 56
$ws.Range("H19").Value = -7.142857142857
$ws.Range("I19").Value = 287
$ws.Range("J19").Value = 280
$ws.Range("K19").Value = 2.5
$ws.Range("L19").Value = 36.666666666666
$ws.Range("M19").Value = 43.5
$ws.Range("N19").Value = 19.583333333333
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 18
$ws.Range("H20").Value = 38.461538461538
$ws.Range("I20").Value = 69
$ws.Range("J20").Value = 48
$ws.Range("K20").Value = 43.75
$ws.Range("L20").Value = 146.428571428571
$ws.Range("M20").Value = 72.5
$ws.Range("N20").Value = -87.126865671641
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 35
$ws.Range("E21").Value = -51.428571428571
$ws.Range("F21").Value = 108
$ws.Range("G21").Value = 106
$ws.Range("H21").Value = 1.886792452830
$ws.Range("I21").Value = 612
$ws.Range("J21").Value = 518
$ws.Range("K21").Value = 18.146718146718
$ws.Range("L21").Value = 51.861042183622
$ws.Range("M21").Value = 46.062052505966
$ws.Range("N21").Value = -63.462686567164
$ws.Range("F22").Value = 1
$ws.Range("M22").Value = -29.411764705882
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = -25
$ws.Range("J23").Value = 10
$ws.Range("K23").Value = 20
$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 26
$ws.Range("E24").Value = -7.692307692307
$ws.Range("F24").Value = 95
$ws.Range("G24").Value = 106
$ws.Range("H24").Value = -10.377358490566
$ws.Range("I24").Value = 768
$ws.Range("J24").Value = 629
$ws.Range("K24").Value = 22.098569157392
$ws.Range("L24").Value = 69.911504424778
$ws.Range("M24").Value = 82.422802850356
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 14.285714285714
$ws.Range("F25").Value = 22
$ws.Range("G25").Value = 32
$ws.Range("H25").Value = -31.25
$ws.Range("I25").Value = 149
$ws.Range("J25").Value = 151
$ws.Range("K25").Value = -1.324503311258
$ws.Range("L25").Value = 27.350427350427
$ws.Range("M25").Value = 36.697247706422
$ws.Range("F26").Value = 1
$ws.Range("H26").Value = 0
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 40
$ws.Range("I27").Value = 24
$ws.Range("J27").Value = 27
$ws.Range("K27").Value = -11.111111111111
$ws.Range("L27").Value = 84.615384615384
$ws.Range("J28").Value = 4
$ws.Range("K28").Value = -75
$ws.Range("J29").Value = 4
$ws.Range("K29").Value = -75

# --- Cells that become "N/A" placeholder text ("0" or "***.*") ---
# Force literal text via a leading apostrophe, then copy the right-aligned
# General-format style (from A15, used throughout for these text placeholders)
# over the cell so only the formatting (not the value) is restored.
$ws.Range("C15").Formula = "'0"
$ws.Range("A15").Copy()
$ws.Range("C15").PasteSpecial($xlPasteFormats)
$ws.Range("D15").Formula = "'0"
$ws.Range("A15").Copy()
$ws.Range("D15").PasteSpecial($xlPasteFormats)
$ws.Range("E15").Formula = "'***.*"
$ws.Range("A15").Copy()
$ws.Range("E15").PasteSpecial($xlPasteFormats)
$ws.Range("C16").Formula = "'0"
$ws.Range("A15").Copy()
$ws.Range("C16").PasteSpecial($xlPasteFormats)
$ws.Range("C18").Formula = "'0"
$ws.Range("A15").Copy()
$ws.Range("C18").PasteSpecial($xlPasteFormats)
$ws.Range("G22").Formula = "'0"
$ws.Range("A15").Copy()
$ws.Range("G22").PasteSpecial($xlPasteFormats)
$ws.Range("H22").Formula = "'***.*"
$ws.Range("A15").Copy()
$ws.Range("H22").PasteSpecial($xlPasteFormats)
$ws.Range("C26").Formula = "'0"
$ws.Range("A15").Copy()
$ws.Range("C26").PasteSpecial($xlPasteFormats)
$ws.Range("D26").Formula = "'0"
$ws.Range("A15").Copy()
$ws.Range("D26").PasteSpecial($xlPasteFormats)
$ws.Range("E26").Formula = "'***.*"
$ws.Range("A15").Copy()
$ws.Range("E26").PasteSpecial($xlPasteFormats)
$ws.Range("C27").Formula = "'0"
$ws.Range("A15").Copy()
$ws.Range("C27").PasteSpecial($xlPasteFormats)

# --- Cells that change from "N/A" placeholder text back to real numbers ---
# Setting a numeric value already clears the text flag; re-copy the donor
# numeric style (count style from I14, percentage style from K14) to match formatting.
$ws.Range("D23").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D23").PasteSpecial($xlPasteFormats)
$ws.Range("E23").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E23").PasteSpecial($xlPasteFormats)
$ws.Range("D27").Value = 3
$ws.Range("I14").Copy()
$ws.Range("D27").PasteSpecial($xlPasteFormats)
$ws.Range("E27").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E27").PasteSpecial($xlPasteFormats)
$ws.Range("D28").Value = 2
$ws.Range("I14").Copy()
$ws.Range("D28").PasteSpecial($xlPasteFormats)
$ws.Range("E28").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E28").PasteSpecial($xlPasteFormats)
$ws.Range("G28").Value = 2
$ws.Range("I14").Copy()
$ws.Range("G28").PasteSpecial($xlPasteFormats)
$ws.Range("H28").Value = -100
$ws.Range("K14").Copy()
$ws.Range("H28").PasteSpecial($xlPasteFormats)
$ws.Range("D29").Value = 2
$ws.Range("I14").Copy()
$ws.Range("D29").PasteSpecial($xlPasteFormats)
$ws.Range("E29").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E29").PasteSpecial($xlPasteFormats)
$ws.Range("G29").Value = 2
$ws.Range("I14").Copy()
$ws.Range("G29").PasteSpecial($xlPasteFormats)
$ws.Range("H29").Value = -100
$ws.Range("K14").Copy()
$ws.Range("H29").PasteSpecial($xlPasteFormats)
